# Apply the edits described by the commit:
#  - Insert a new catalog item "Luminárias de móvel USB" (Caixa 35 / Iluminação /
#    Quarto do Bento / Média priority) right before the current row 263
#    ("Caixa 36" / Vestuário / Vestido H&M ...), shifting all subsequent rows down by one.
#  - Clear the (already blank) priority-note cells F196:F202 so they are stored as
#    blank numeric cells instead of blank inline-string cells.
#  - Update the "Resumo" sheet totals: Total de Itens 437 -> 438, Média 272 -> 273.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Catálogo de Mudança")

# --- Insert new row 263 ------------------------------------------------
$ws.Rows.Item(263).Insert()

# Copy cell formatting (borders / wrap / alignment / fill / font) from the row that
# just got pushed down to 264 so the new row looks like every other data row.
$srcRow = $ws.Range("A264:F264")
$dstRow = $ws.Range("A263:F263")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats

# The priority column (E) on this new row is "Média", so copy that specific
# highlight style from an existing "Média" row (row 260) onto E263.
$srcE = $ws.Range("E260")
$dstE = $ws.Range("E263")
$srcE.Copy()
$dstE.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Fill in the values for the newly inserted row.
$ws.Cells.Item(263, 1).Value = "Caixa 35"
$ws.Cells.Item(263, 2).Value = "Iluminação"
$ws.Cells.Item(263, 3).Value = "Luminárias de móvel USB"
$ws.Cells.Item(263, 4).Value = "Quarto do Bento"
$ws.Cells.Item(263, 5).Value = "Média"
$ws.Cells.Item(263, 6).ClearContents()

# --- Clear the blank F196:F202 notes cells -----------------------------
$ws.Range("F196:F202").ClearContents()

# --- Update the Resumo sheet summary numbers ---------------------------
$ws2 = $wb.Worksheets.Item("Resumo")
$ws2.Cells.Item(3, 2).Value = 438
$ws2.Cells.Item(7, 2).Value = 273
